# Insert a new weekly price record as the first observation (row 858) for
# "Hortaliza, Femacal de La Calera - Zapallo italiano", pushing the existing
# rows 858:877 down to 859:878.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 858; this shifts rows
# 858:877 down to 859:878, preserving all of their data and formatting.
$ws.Rows.Item(858).Insert()

# Populate the newly inserted row 858 with the new observation.
$ws.Cells.Item(858, 1).Value = 3
$ws.Cells.Item(858, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(858, 3).Value = "Coquimbo"
$ws.Cells.Item(858, 4).Value = 45239
$ws.Cells.Item(858, 5).Value = 5
$ws.Cells.Item(858, 6).Value = 100112032
$ws.Cells.Item(858, 7).Value = "Zapallo italiano"
$ws.Cells.Item(858, 8).Value = "Sin especificar"
$ws.Cells.Item(858, 9).Value = "Primera"
$ws.Cells.Item(858, 10).Value = 65
$ws.Cells.Item(858, 11).Value = 9000
$ws.Cells.Item(858, 12).Value = 9000
$ws.Cells.Item(858, 13).Value = 9000
$ws.Cells.Item(858, 14).Value = "$/caja 36 unidades"
$ws.Cells.Item(858, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(858, 16).Value = 250
$ws.Cells.Item(858, 17).Value = 36
$ws.Cells.Item(858, 18).Value = "Hortaliza"
